$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    4  = @{ C = 5.907218141265402;    E = 5.917486466529609 }
    5  = @{ C = 6.511263427347003;    E = 6.422943767670297 }
    6  = @{ C = 4.325828829470257;    E = 4.493586323244281 }
    7  = @{ C = 2.964652118442834;    E = 3.83627393798931 }
    8  = @{ C = 3.171852776411788;    E = 4.030605385534614 }
    9  = @{ C = 1.773712379859993;    E = 3.484530515673856 }
    10 = @{ C = 2.533350906619081;    E = 3.524103740130435 }
    11 = @{ C = 2.661040979345697;    E = 3.567108445582057 }
    12 = @{ C = 3.150198973767537;    E = 3.699072253610103 }
    13 = @{ C = 0.4641929091049102;   E = 2.550259844884462 }
    14 = @{ C = 2.585454129751663;    E = 2.671828487424377 }
    15 = @{ C = -0.4532848472497908;  E = 2.066462658785673 }
    16 = @{ C = 0.9477102747197819;   E = 1.83067479293082 }
    17 = @{ C = 2.222852754198135;    E = 1.934107558751452 }
    18 = @{ C = -0.007094633234694392; E = 1.444584248586422 }
    19 = @{ C = 3.078872076370009;    E = 2.279508996785351 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
